$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": refresh the handoff-status report for the
# file that was just handed off (0098402b-a623-465f-9c4e-0baa758574cd.md),
# row 4 on every sheet.

# Overview sheet: Latest Handoff Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D4").Value = "2016-03-21 03:03:18"

# zh-cn sheet: Latest Handoff Datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-21 03:03:09"

# de-de sheet: Latest Handoff Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-21 03:03:18"
